$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Solver scenario results: fill in the "Variables" allocation table (B35:D37)
# with the solution the Solver add-in produced, then wire up the Objective
# and Constraints formulas that reference it.
# ---------------------------------------------------------------------------
$ws.Range("B35").Value = 40
$ws.Range("C35").Value = 40
$ws.Range("D35").Value = 80
$ws.Range("B36").Value = 100
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 0
$ws.Range("B37").Value = 0
$ws.Range("C37").Value = 40
$ws.Range("D37").Value = 0

$ws.Range("B40").Formula = "=SUMPRODUCT(B35:D37,B17:D19)"

$ws.Range("B45").Formula = "=SUMPRODUCT(B35:D35,B17:D17)"
$ws.Range("B46").Formula = "=SUMPRODUCT(B36:D36,B18:D18)"
$ws.Range("B47").Formula = "=SUMPRODUCT(B37:D37,B19:D19)"

$ws.Range("B50").Formula = "=SUM(B35:B37)"
$ws.Range("B51").Formula = "=SUM(C35:C37)"
$ws.Range("B52").Formula = "=SUM(D35:D37)"

# ---------------------------------------------------------------------------
# Solver parameters saved with the sheet (Data > Solver dialog settings).
# These are stored as hidden, sheet-scoped defined names named "solver_*".
# ---------------------------------------------------------------------------
function Add-HiddenName($name, $refersTo) {
    $nm = $ws.Names.Add($name, $refersTo)
    $nm.Visible = $false
}

Add-HiddenName "solver_adj"  '=Sheet1!$B$35:$D$37'
Add-HiddenName "solver_cvg"  '0.0001'
Add-HiddenName "solver_drv"  '1'
Add-HiddenName "solver_eng"  '2'
Add-HiddenName "solver_itr"  '2147483647'
Add-HiddenName "solver_lhs1" '=Sheet1!$B$45:$B$47'
Add-HiddenName "solver_lhs2" '=Sheet1!$B$50:$B$52'
Add-HiddenName "solver_lin"  '1'
Add-HiddenName "solver_mip"  '2147483647'
Add-HiddenName "solver_mni"  '30'
Add-HiddenName "solver_mrt"  '0.075'
Add-HiddenName "solver_msl"  '2'
Add-HiddenName "solver_neg"  '1'
Add-HiddenName "solver_nod"  '2147483647'
Add-HiddenName "solver_num"  '2'
Add-HiddenName "solver_opt"  '=Sheet1!$B$40'
Add-HiddenName "solver_pre"  '0.000001'
Add-HiddenName "solver_rbv"  '1'
Add-HiddenName "solver_rel1" '1'
Add-HiddenName "solver_rel2" '1'
Add-HiddenName "solver_rhs1" '=Sheet1!$D$45:$D$47'
Add-HiddenName "solver_rhs2" '=Sheet1!$D$50:$D$52'
Add-HiddenName "solver_rlx"  '2'
Add-HiddenName "solver_rsd"  '0'
Add-HiddenName "solver_scl"  '1'
Add-HiddenName "solver_sho"  '2'
Add-HiddenName "solver_ssz"  '100'
Add-HiddenName "solver_tim"  '2147483647'
Add-HiddenName "solver_tol"  '0.01'
Add-HiddenName "solver_typ"  '1'
Add-HiddenName "solver_val"  '0'
Add-HiddenName "solver_ver"  '2'

# ---------------------------------------------------------------------------
# View state: the sheet was left scrolled near the bottom of the Solver
# results with B47 selected.
# ---------------------------------------------------------------------------
$ws.Range("B47").Select()
